# Weekly fruit/vegetable data refresh: a new week's price record was
# inserted at the top of this subset's data block (row 247), pushing the
# existing records (old rows 247-329) down by one row to 248-330.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 247, shifting rows 247:329 down to 248:330.
$ws.Rows("247").Insert()

# Populate the newly inserted row 247 with this week's record.
$ws.Range("A247").Value = 4
$ws.Range("B247").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C247").Value = "Los Lagos"
$ws.Range("D247").Value = 44809
$ws.Range("E247").Value = 10
$ws.Range("F247").Value = 100112037
$ws.Range("G247").Value = "Cebollín"
$ws.Range("H247").Value = "Sin especificar"
$ws.Range("I247").Value = "Primera"
$ws.Range("J247").Value = 70
$ws.Range("K247").Value = 11000
$ws.Range("L247").Value = 11000
$ws.Range("M247").Value = 11000
$ws.Range("N247").Value = "`$/paquete 36 unidades"
$ws.Range("O247").Value = "Región Metropolitana"
$ws.Range("P247").Value = 306
$ws.Range("Q247").Value = 36
$ws.Range("R247").Value = "Hortaliza"
